$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H32").Value = 3132.3635
$ws.Range("J32").Value = 3162.4443
$ws.Range("L32").Value = 3162.4443
$ws.Range("N32").Value = -3814.4443
$ws.Range("H33").Value = 448.31708
$ws.Range("J33").Value = 1284.125
$ws.Range("L33").Value = 1284.125
$ws.Range("N33").Value = -1742.125
$ws.Range("H39").Value = 287.18182
$ws.Range("I39").Value = 183.33333
$ws.Range("J39").Value = 754.5
$ws.Range("K39").Value = 549.99999
$ws.Range("L39").Value = 2263.5
$ws.Range("M39").Value = -253.99999
$ws.Range("N39").Value = -2855.5
$ws.Range("H40").Value = 3496.5334
$ws.Range("J40").Value = 3777.182
$ws.Range("L40").Value = 3777.182
$ws.Range("N40").Value = -4127.182
$ws.Range("H70").Value = 35823284
$ws.Range("I70").Value = 250995
$ws.Range("J70").Value = 62502500
$ws.Range("K70").Value = 752985
$ws.Range("L70").Value = 187507500
$ws.Range("M70").Value = -752715
$ws.Range("N70").Value = -187508040
$ws.Range("H73").Value = 35823284
$ws.Range("I73").Value = 250995
$ws.Range("J73").Value = 62502500
$ws.Range("K73").Value = 752985
$ws.Range("L73").Value = 187507500
$ws.Range("M73").Value = -752049
$ws.Range("N73").Value = -187509372
$ws.Range("H116").Value = 50029436
$ws.Range("I116").Value = 55587150
$ws.Range("K116").Value = 55587150
$ws.Range("M116").Value = -55583708
$ws.Range("H118").Value = 1809.75
$ws.Range("I118").Value = 1740
$ws.Range("K118").Value = 5220
$ws.Range("M118").Value = -3563
$ws.Range("H138").Value = 7586.3135
$ws.Range("I138").Value = 5472.75
$ws.Range("J138").Value = 7872.8984
$ws.Range("K138").Value = 16418.25
$ws.Range("L138").Value = 23618.6952
$ws.Range("M138").Value = -11278.25
$ws.Range("N138").Value = -33898.6952

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H74").Value = 2510.2942
$ws.Range("I74").Value = 1543.5476
$ws.Range("J74").Value = 7021.778
$ws.Range("K74").Value = 1543.5476
$ws.Range("L74").Value = 7021.778
$ws.Range("M74").Value = -669.5476000000001
$ws.Range("N74").Value = -8769.778
$ws.Range("H77").Value = 2510.2942
$ws.Range("I77").Value = 1543.5476
$ws.Range("J77").Value = 7021.778
$ws.Range("K77").Value = 7717.738
$ws.Range("L77").Value = 35108.89
$ws.Range("M77").Value = -3349.738
$ws.Range("N77").Value = -43844.89
$ws.Range("H114").Value = 47493.5
$ws.Range("J114").Value = 47493.5
$ws.Range("L114").Value = 47493.5
$ws.Range("N114").Value = -56171.5
$ws.Range("H119").Value = 99972
$ws.Range("J119").Value = 99972
$ws.Range("L119").Value = 99972
$ws.Range("N119").Value = -109648
$ws.Range("H122").Value = 7178.6597
$ws.Range("I122").Value = 6900.2583
$ws.Range("J122").Value = 7718.0625
$ws.Range("K122").Value = 20700.7749
$ws.Range("L122").Value = 23154.1875
$ws.Range("M122").Value = -18250.7749
$ws.Range("N122").Value = -28054.1875
$ws.Range("H132").Value = 325582.53
$ws.Range("I132").Value = 456730.12
$ws.Range("K132").Value = 1370190.36
$ws.Range("M132").Value = -1367660.36

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H29").Value = 1503.25
$ws.Range("I29").Value = 1503.25
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 1503.25
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -1214.25

$ws.Range("N29").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H69").Value = 37624.57
$ws.Range("I69").Value = 18243.25
$ws.Range("K69").Value = 18243.25
$ws.Range("M69").Value = -17494.25
$ws.Range("H72").Value = 37624.57
$ws.Range("I72").Value = 18243.25
$ws.Range("K72").Value = 54729.75
$ws.Range("M72").Value = -50985.75
$ws.Range("H99").Value = 2674.3333
$ws.Range("J99").Value = 2498.3333
$ws.Range("L99").Value = 2498.3333
$ws.Range("N99").Value = -5494.3333
$ws.Range("H126").Value = 2674.3333
$ws.Range("J126").Value = 2498.3333
$ws.Range("L126").Value = 7494.999899999999
$ws.Range("N126").Value = -12434.9999
$ws.Range("H132").Value = 2918.0557
$ws.Range("I132").Value = 2607.8125
$ws.Range("K132").Value = 7823.4375
$ws.Range("M132").Value = -5293.4375

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H12").Value = 98.625
$ws.Range("J12").Value = 118.8
$ws.Range("L12").Value = 356.4
$ws.Range("N12").Value = -702.4
$ws.Range("H16").Value = 600.25
$ws.Range("I16").Value = 600.5
$ws.Range("J16").Value = 600
$ws.Range("K16").Value = 1801.5
$ws.Range("L16").Value = 1800
$ws.Range("M16").Value = -1628.5
$ws.Range("N16").Value = -2146
$ws.Range("H32").Value = 599
$ws.Range("J32").Value = 198
$ws.Range("L32").Value = 594
$ws.Range("N32").Value = -1160

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("H102").Value = 2199.8928
$ws.Range("I102").Value = 1994.6666
$ws.Range("K102").Value = 1994.6666
$ws.Range("M102").Value = -372.6666
$ws.Range("H107").Value = 316.66666
$ws.Range("I107").Value = 325
$ws.Range("J107").Value = 300
$ws.Range("K107").Value = 325
$ws.Range("L107").Value = 300
$ws.Range("M107").Value = 1595
$ws.Range("N107").Value = -4140
$ws.Range("H121").Value = 115000
$ws.Range("J121").Value = 115000
$ws.Range("L121").Value = 115000
$ws.Range("N121").Value = -118494

$ws.Range("N75").ClearContents()
$ws.Range("N78").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H7").Value = 7837.6665
$ws.Range("J7").Value = 10941.429
$ws.Range("L7").Value = 10941.429
$ws.Range("N7").Value = -11165.429
$ws.Range("H22").Value = 1247.25
$ws.Range("I22").Value = 1500
$ws.Range("J22").Value = 994.5
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 994.5
$ws.Range("M22").Value = -1205
$ws.Range("N22").Value = -1584.5
$ws.Range("H27").Value = 1247.25
$ws.Range("I27").Value = 1500
$ws.Range("J27").Value = 994.5
$ws.Range("K27").Value = 1500
$ws.Range("L27").Value = 994.5
$ws.Range("M27").Value = -1393
$ws.Range("N27").Value = -1208.5
$ws.Range("H104").Value = 85998.336
$ws.Range("J104").Value = 85998.336
$ws.Range("L104").Value = 85998.336
$ws.Range("N104").Value = -92986.336
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("H126").Value = 7837.6665
$ws.Range("J126").Value = 10941.429
$ws.Range("L126").Value = 32824.287
$ws.Range("N126").Value = -37764.287
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("H136").Value = 8971
$ws.Range("J136").Value = 10184.75
$ws.Range("L136").Value = 30554.25
$ws.Range("N136").Value = -35654.25

$ws.Range("N119").ClearContents()
$ws.Range("N129").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H54").Value = 22000
$ws.Range("J54").Value = 26000
$ws.Range("L54").Value = 26000
$ws.Range("N54").Value = -27040
$ws.Range("H70").Value = 43124.812
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("H73").Value = 43124.812
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("H107").Value = 526.8889
$ws.Range("J107").Value = 544.9
$ws.Range("L107").Value = 1634.7
$ws.Range("N107").Value = -5474.7
$ws.Range("H119").Value = 97494.5
$ws.Range("J119").Value = 97494.5
$ws.Range("L119").Value = 97494.5
$ws.Range("N119").Value = -107170.5
$ws.Range("H122").Value = 32260422
$ws.Range("I122").Value = 55557424
$ws.Range("K122").Value = 166672272
$ws.Range("M122").Value = -166669822
$ws.Range("H129").Value = 120000.664
$ws.Range("J129").Value = 120000.664
$ws.Range("L129").Value = 120000.664
$ws.Range("N129").Value = -130000.664

$ws.Range("M70").ClearContents()
$ws.Range("M73").ClearContents()
